$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "97.265.98"
$ws.Range("E2").Value = "  +2.83%  "
$ws.Range("D3").Value = "3.351.90"
$ws.Range("E3").Value = "  +7.63%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.13"
$ws.Range("E5").Value = "  +1.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "623.30"
$ws.Range("E6").Value = "  +1.22%  "
$ws.Range("E7").Value = "  -0.47%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.387"
$ws.Range("E8").Value = "  -1.46%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  +0.13%  "
$ws.Range("D10").Value = "3.349.42"
$ws.Range("E10").Value = "  +7.52%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.788"
$ws.Range("E11").Value = "  -4.86%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.199"
$ws.Range("E12").Value = "  +0.39%  "
$ws.Range("D13").Value = "97.122.81"
$ws.Range("E13").Value = "  +3.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000246"
$ws.Range("E14").Value = "  +0.30%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "3.975.09"
$ws.Range("E15").Value = "  +7.61%  "
$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "35.02"
$ws.Range("E16").Value = "  +0.93%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.50"
$ws.Range("E17").Value = "  +1.74%  "
$ws.Range("D18").Value = "3.338.63"
$ws.Range("E18").Value = "  +6.94%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.56"
$ws.Range("E19").Value = "  -4.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.12"
$ws.Range("E20").Value = "  +0.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "488.05"
$ws.Range("E21").Value = "  +8.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0000209"
$ws.Range("E22").Value = "  +4.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.83"
$ws.Range("E23").Value = "  -2.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.28"
$ws.Range("E24").Value = "  +2.94%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.67"
$ws.Range("E25").Value = "  +0.35%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "88.35"
$ws.Range("E26").Value = "  +2.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.08"
$ws.Range("E27").Value = "  +0.34%  "
$ws.Range("D28").Value = "3.537.60"
$ws.Range("E28").Value = "  +7.62%  "
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("B30").Value = "Cronos"
$ws.Range("C30").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.181"
$ws.Range("E30").Value = "  +0.71%  "
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.242"
$ws.Range("E31").Value = "  -6.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.999"
$ws.Range("E32").Value = "  -0.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.120"
$ws.Range("E33").Value = "  -5.63%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "9.29"
$ws.Range("E34").Value = "  -0.61%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "27.64"
$ws.Range("E35").Value = "  +5.35%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.40"
$ws.Range("E36").Value = "  -6.32%  "
$ws.Range("E37").Value = "  -6.68%  "
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "501.93"
$ws.Range("E38").Value = "  +4.53%  "
$ws.Range("B39").Value = "PancakeSwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.93"
$ws.Range("E39").Value = "  +0.63%  "
$ws.Range("E40").Value = "  +2.85%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.449"
$ws.Range("E41").Value = "  -1.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.27"
$ws.Range("E42").Value = "  -1.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.805"
$ws.Range("E43").Value = "  +16.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.26"
$ws.Range("E44").Value = "  -1.54%  "
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("E46").Value = "  -7.65%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "159.57"
$ws.Range("E47").Value = "  -0.60%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.92"
$ws.Range("E48").Value = "  +2.99%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.54"
$ws.Range("E49").Value = "  +2.37%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "45.19"
$ws.Range("E50").Value = "  +2.95%  "
$ws.Range("E51").Value = "  +1.70%  "
